$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new example paragraph:
#        data.add("https://www.google.com/");
#    right before the "sheet.getData().add(data);" paragraph (i.e. right
#    after the "data.add("4.1");" paragraph), matching the indentation
#    (left=720, no "space after") of its sibling example paragraphs.
# ---------------------------------------------------------------------

$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like 'data.add("4.1")*') {
        $anchor = $p
    }
}

$anchor.Range.InsertParagraphAfter()

# Locate the freshly-minted (empty) paragraph that now sits between
# "data.add("4.1");" and "sheet.getData().add(data);".
$newPara = $anchor.Next()
$pos = $newPara.Range.Start

foreach ($chunk in @("data.add", "(", '"', "https://www.google.com/", '"', ");")) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($chunk)
    $pos = $pos + $chunk.Length
}

# ---------------------------------------------------------------------
# 2. Expand the "If the value set for a cell appears to be an integer
#    or floating-point number ..." paragraph with a new sentence about
#    URL values being turned into hyperlinks, and relocate the
#    "_GoBack" bookmark to the end of the paragraph (after all text).
# ---------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*integer or*") {
        $target = $p
    }
}

$paraEnd = $target.Range.End
$contentEnd = $paraEnd - 1    # position right before the paragraph mark

$newSentence = "If the value for the cell appears to be a URL, then the EEH will set the text as a clickable hyperlink in the cell. "

# NOTE: adding a zero-length bookmark whose position sits exactly one
# character before a paragraph mark confuses this engine's Bookmarks.Add
# (it resolves to some unrelated {0,20} range). Work around it by
# temporarily inserting a one-character placeholder right before the
# paragraph mark, doing all the real text/bookmark surgery "before" that
# placeholder (so the bookmark position is never adjacent to the CR),
# and deleting the placeholder afterwards.
$ph = $d.Range($contentEnd, $contentEnd)
$ph.InsertAfter("Z")

$ins2 = $d.Range($contentEnd, $contentEnd)
$ins2.InsertAfter($newSentence)

$newContentEnd = $contentEnd + $newSentence.Length

$d.Bookmarks.Item("_GoBack").Delete()
$bmRange = $d.Range($newContentEnd, $newContentEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character again.
$phRange = $d.Range($newContentEnd, $newContentEnd + 1)
$phRange.Delete()

Write-Output "done"
